# Adds two new Mac-Address user rows (rows 31 and 32) to Sheet1, mirroring
# the values of the existing rows above them.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 32 (John Doe) is populated first so its text lands earlier in the
# workbook's shared-string table than row 31's (Jane Smith), matching the
# order the strings were originally authored in.
$ws.Cells.Item(32, 1).Value = 110031
$ws.Cells.Item(32, 2).Value = 9317596767
$ws.Cells.Item(32, 3).Value = "John Doe"
$ws.Cells.Item(32, 4).Value = "john.doe@xyz.com"
$ws.Cells.Item(32, 5).Value = 818876431
$ws.Cells.Item(32, 6).Value = "ACT"
$ws.Cells.Item(32, 7).Value = "eng"
$ws.Cells.Item(32, 8).Value = "PWD"
$ws.Cells.Item(32, 9).Value = $true
$ws.Cells.Item(32, 9).HorizontalAlignment = -4131
$ws.Cells.Item(32, 10).Value = "superadmin"
$ws.Cells.Item(32, 11).Value = "now()"
$ws.Cells.Item(32, 12).Value = "now()"

# Row 31: Jane Smith
$ws.Cells.Item(31, 1).Value = 110030
$ws.Cells.Item(31, 2).Value = 9317596768
$ws.Cells.Item(31, 3).Value = "Jane Smith"
$ws.Cells.Item(31, 4).Value = "jane.smith@xyz.com"
$ws.Cells.Item(31, 5).Value = 818876432
$ws.Cells.Item(31, 6).Value = "ACT"
$ws.Cells.Item(31, 7).Value = "eng"
$ws.Cells.Item(31, 8).Value = "PWD"
$ws.Cells.Item(31, 9).Value = $true
$ws.Cells.Item(31, 9).HorizontalAlignment = -4131
$ws.Cells.Item(31, 10).Value = "superadmin"
$ws.Cells.Item(31, 11).Value = "now()"
$ws.Cells.Item(31, 12).Value = "now()"

# Update the selection to match the post-edit state recorded in the workbook
$ws.Range("F30").Select()
